$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(9, 1).Value = 123853472
$ws.Cells.Item(9, 2).Value = 123
$ws.Cells.Item(9, 3).Value = 1

$ws.Cells.Item(10, 1).Value = 124578963
$ws.Cells.Item(10, 2).Value = 123
$ws.Cells.Item(10, 3).Value = 1
